$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the sheet name in the workbook
$ws.Name = "Through 2021-10-06"

# Row 12 (October) updates
$ws.Range("A12").Value = "October (through 10-06)"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 0.1667
$ws.Range("F12").Value = 9
$ws.Range("L12").Value = 15
$ws.Range("O12").Value = 7
$ws.Range("R12").Value = 33
$ws.Range("U12").Value = 46

# Row 13 (Total) updates
$ws.Range("C13").Value = 201
$ws.Range("D13").Value = 0.1336
$ws.Range("F13").Value = 392
$ws.Range("G13").Value = 0.105
$ws.Range("L13").Value = 502
$ws.Range("M13").Value = 0.1083
$ws.Range("O13").Value = 386
$ws.Range("P13").Value = 0.1002
$ws.Range("R13").Value = 881
$ws.Range("S13").Value = 0.0567
$ws.Range("U13").Value = 1216
$ws.Range("V13").Value = 0.061
